$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 215, pushing all existing rows (215-275) down to 217-277.
$ws.Rows("215:216").Insert()

# New row 215: weekly price entry for Murcott / Primera
$ws.Cells.Item(215, 1).Value = 11
$ws.Cells.Item(215, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(215, 3).Value = "Bíobío"
$ws.Cells.Item(215, 4).Value = 45204
$ws.Cells.Item(215, 5).Value = 8
$ws.Cells.Item(215, 6).Value = "Fruta"
$ws.Cells.Item(215, 7).Value = 100102
$ws.Cells.Item(215, 8).Value = "Cítricos"
$ws.Cells.Item(215, 9).Value = 100102004
$ws.Cells.Item(215, 10).Value = "Mandarina"
$ws.Cells.Item(215, 11).Value = "Murcott"
$ws.Cells.Item(215, 12).Value = "Primera"
$ws.Cells.Item(215, 13).Value = 270
$ws.Cells.Item(215, 14).Value = 8000
$ws.Cells.Item(215, 15).Value = 8500
$ws.Cells.Item(215, 16).Value = 8222
$ws.Cells.Item(215, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(215, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(215, 19).Value = 457
$ws.Cells.Item(215, 20).Value = 18

# New row 216: weekly price entry for Murcott / Segunda
$ws.Cells.Item(216, 1).Value = 11
$ws.Cells.Item(216, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(216, 3).Value = "Bíobío"
$ws.Cells.Item(216, 4).Value = 45204
$ws.Cells.Item(216, 5).Value = 8
$ws.Cells.Item(216, 6).Value = "Fruta"
$ws.Cells.Item(216, 7).Value = 100102
$ws.Cells.Item(216, 8).Value = "Cítricos"
$ws.Cells.Item(216, 9).Value = 100102004
$ws.Cells.Item(216, 10).Value = "Mandarina"
$ws.Cells.Item(216, 11).Value = "Murcott"
$ws.Cells.Item(216, 12).Value = "Segunda"
$ws.Cells.Item(216, 13).Value = 200
$ws.Cells.Item(216, 14).Value = 7000
$ws.Cells.Item(216, 15).Value = 7000
$ws.Cells.Item(216, 16).Value = 7000
$ws.Cells.Item(216, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(216, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(216, 19).Value = 389
$ws.Cells.Item(216, 20).Value = 18
